$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.622.83"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "1.643.74"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'215.83"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").Value = "1.872.36"
$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("D13").Value = "'4.21"
$ws.Range("E13").Value = "  +3.16%  "

$ws.Range("D14").Value = "1.638.06"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").Value = "'65.96"
$ws.Range("E16").Value = "  +4.32%  "

$ws.Range("D17").Value = "26.666.84"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "'218.57"
$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("E22").Value = "  +2.12%  "

$ws.Range("D23").Value = "'9.56"
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("E24").Value = "  +11.26%  "

$ws.Range("D25").Value = "'146.24"
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("E29").Value = "  +2.46%  "

$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "  +2.86%  "

$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("E33").Value = "  +2.51%  "

$ws.Range("D34").Value = "1.274.63"
$ws.Range("E34").Value = "  +5.38%  "

$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("E36").Value = "  +6.25%  "

$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("E38").Value = "  +5.82%  "

$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  +2.53%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = "  +2.41%  "

$ws.Range("E42").Value = "  -1.57%  "

$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("D44").Value = "1.783.72"
$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("D45").Value = "'93.16"
$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("D46").Value = "'59.85"
$ws.Range("E46").Value = "  +9.39%  "

$ws.Range("E47").Value = "  +3.34%  "

$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").Value = "'7.78"
$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("E51").Value = "  -0.60%  "
